# Add 2022-Q1 data.
#
# The existing "总计" (Total) summary sheet is renamed to "2022-Q1" and its
# body is replaced with the quarter's fund-holding detail table (columns
# A-H). A brand new "总计" sheet is appended at the end, carrying the old
# summary table plus a new first data row for "2022-Q1".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Repurpose the current "总计" sheet as the new "2022-Q1" detail sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Clear out the old 4-column summary body before writing the new table.
$q1.Range("A1:D6").ClearContents()

# Grab a cell that already carries the bold/centered/bordered header style
# (style index 2 in the source workbook) so we can stamp it onto the new
# header row and the row-index column without inventing new style slots.
$styleSrc = $wb.Worksheets.Item("2021-Q4").Range("B1")

# Header row.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$styleSrc.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# Fund rows. Columns B, D, E, F, G hold numeric-looking text (fund codes
# with leading zeros, percentages, NAV figures) that must stay text, so the
# range is pre-formatted as "@" (Text) before the values land.
$q1.Range("B2:G9").NumberFormat = "@"

$funds = @(
    @("166027", "中欧创业板两年定期开放混合A",     "21.11", "99.64", "5.26", "1.1104", 5),
    @("160726", "嘉实瑞享定期开放灵活配置混合",     "23.58", "63.95", "1.99", "0.4692", 10),
    @("009138", "嘉实瑞成两年持有期混合A",         "22.27", "75.95", "2.10", "0.4677", 9),
    @("009791", "中欧创业板两年定期开放混合C",     "5.21",  "99.64", "5.26", "0.2740", 5),
    @("009139", "嘉实瑞成两年持有期混合C",         "4.37",  "75.95", "2.10", "0.0918", 9),
    @("164826", "工银瑞信创业板两年定期开放混合A", "2.02",  "72.79", "3.75", "0.0758", 7),
    @("010889", "工银瑞信创业板两年定期开放混合C", "0.16",  "72.79", "3.75", "0.0060", 7),
    @("005167", "嘉实润泽量化一年定期开放混合",     "0.56",  "27.26", "0.66", "0.0037", 7)
)

for ($i = 0; $i -lt $funds.Length; $i++) {
    $row = $i + 2
    $f = $funds[$i]
    $q1.Cells.Item($row, 1).Value = $i
    $q1.Range("B$row").Value = $f[0]
    $q1.Range("C$row").Value = $f[1]
    $q1.Range("D$row").Value = $f[2]
    $q1.Range("E$row").Value = $f[3]
    $q1.Range("F$row").Value = $f[4]
    $q1.Range("G$row").Value = $f[5]
    $q1.Cells.Item($row, 8).Value = $f[6]
}

$styleSrc.Copy()
$q1.Range("A2:A9").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Append a brand new "总计" sheet with the updated summary table.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add($null, $lastSheet)
$total.Name = "总计"

# Match the page margins used throughout the rest of the workbook.
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36
$total.Outline.SummaryRow = 1
$total.Outline.SummaryColumn = 1

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$styleSrc.Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$summary = @(
    @("2022-Q1", 8, 2.5),
    @("2021-Q4", 4, 1.52),
    @("2021-Q3", 3, 1.4),
    @("2021-Q2", 3, 1.72),
    @("2021-Q1", 7, 1.03),
    @("2020-Q4", 6, 3.16)
)

for ($i = 0; $i -lt $summary.Length; $i++) {
    $row = $i + 2
    $s = $summary[$i]
    $total.Cells.Item($row, 1).Value = $i
    $total.Range("B$row").Value = $s[0]
    $total.Range("C$row").Value = $s[1]
    $total.Range("D$row").Value = $s[2]
}

$styleSrc.Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
